# Update "Arribos 141" schedule data on sheets TODOS and COMBINADAS (rows 14,
# 22-24, 45-46) and on sheet 215 (row 7), per the commit's refreshed arrival
# times ("Arribos 141 actualizados - 40").

$wb = $excel.ActiveWorkbook

# ---- TODOS ----
$ws = $wb.Worksheets.Item("TODOS")

# Row 14: SCRAPING time shifts 23:18 -> 23:20, MIN 28 -> 26.
$ws.Range("A14").Value = "23:20"
$ws.Range("D14").Value = 26

# Rows 22-24 rotate: the 23:18/23:58 "BANDERA" entries move down one row
# (ETCHEVERRY -> ALUAR -> ABASTO), with the new first row stamped 23:20.
$ws.Range("A22").Value = "23:20"
$ws.Range("B22").Value = "23:58"
$ws.Range("C22").Value = "11X44_ETCHEVERRY"
$ws.Range("D22").Value = 38

$ws.Range("A23").Value = "23:20"
$ws.Range("B23").Value = "23:58"
$ws.Range("C23").Value = "215_ALUAR"
$ws.Range("D23").Value = 38

$ws.Range("A24").Value = ""
$ws.Range("B24").Value = "21:22"
$ws.Range("C24").Value = "15_ABASTO"
$ws.Range("D24").Value = 38

# Rows 45-46 swap: 15_ABASTO moves to row 46, 16_SANTA ANA (now stamped
# 23:20) takes row 45.
$ws.Range("A45").Value = "23:20"
$ws.Range("B45").Value = "00:49"
$ws.Range("C45").Value = "16_SANTA ANA"
$ws.Range("D45").Value = 89

$ws.Range("A46").Value = ""
$ws.Range("B46").Value = "21:23"
$ws.Range("C46").Value = "15_ABASTO"
$ws.Range("D46").Value = 90

# ---- COMBINADAS (mirrors TODOS for this route) ----
$ws = $wb.Worksheets.Item("COMBINADAS")

$ws.Range("A14").Value = "23:20"
$ws.Range("D14").Value = 26

$ws.Range("A22").Value = "23:20"
$ws.Range("B22").Value = "23:58"
$ws.Range("C22").Value = "11X44_ETCHEVERRY"
$ws.Range("D22").Value = 38

$ws.Range("A23").Value = "23:20"
$ws.Range("B23").Value = "23:58"
$ws.Range("C23").Value = "215_ALUAR"
$ws.Range("D23").Value = 38

$ws.Range("A24").Value = ""
$ws.Range("B24").Value = "21:22"
$ws.Range("C24").Value = "15_ABASTO"
$ws.Range("D24").Value = 38

$ws.Range("A45").Value = "23:20"
$ws.Range("B45").Value = "00:49"
$ws.Range("C45").Value = "16_SANTA ANA"
$ws.Range("D45").Value = 89

$ws.Range("A46").Value = ""
$ws.Range("B46").Value = "21:23"
$ws.Range("C46").Value = "15_ABASTO"
$ws.Range("D46").Value = 90

# ---- 215 sheet: single-row refresh (equivalent to old TODOS row 23) ----
$ws = $wb.Worksheets.Item("215")
$ws.Range("A7").Value = "23:20"
$ws.Range("D7").Value = 38
